# DaySale report update:
#  - A new low-stock item ("شفاط ثدي فلاي بيبي") appears in the shortage list, sorted
#    alphabetically just before "كالونا" -> insert a new row at row 32 for it.
#  - The grand-total (column P sum row) grows by the new item's selling price.
#  - The generated-at timestamp in the footer moves from 12:26 PM to 12:32 PM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a fresh row right above the old row 32 ("كالونا ...") - this shifts
#    every row from 32 downward (data rows, the totals row and the footer row)
#    down by one, carrying their merged ranges along with them.
# ---------------------------------------------------------------------------
$ws.Rows("32:32").Insert()

# ---------------------------------------------------------------------------
# 2) Clone the row-32-shaped formatting (borders/fills/fonts/number formats)
#    from the row right below (row 33, which now holds the old "كالونا" row)
#    so the new row matches the rest of the table exactly.
# ---------------------------------------------------------------------------
$ws.Range("A33:Q33").Copy()
$ws.Range("A32:Q32").PasteSpecial(-4122) # xlPasteFormats

# ---------------------------------------------------------------------------
# 3) Fill in the new item's data.
#    Columns L (order limit) and P (selling price) are styled with numeric
#    display formats even though the workbook stores them as plain text, so
#    flip to a text format while assigning, then restore the original number
#    format (keeps the same style id, avoids Excel silently turning the
#    string into a real number / adding a quote-prefix style variant).
# ---------------------------------------------------------------------------
$ws.Range("A32").Value = 26
$ws.Range("C32").Value = "شفاط ثدي فلاي بيبي"
$ws.Range("H32").Value = "0:0"

$limitFmt = $ws.Range("L32").NumberFormat
$ws.Range("L32").NumberFormat = "@"
$ws.Range("L32").Value = "0"
$ws.Range("L32").NumberFormat = $limitFmt

$ws.Range("N32").Value = "65.00"

$priceFmt = $ws.Range("P32").NumberFormat
$ws.Range("P32").NumberFormat = "@"
$ws.Range("P32").Value = "65.0000"
$ws.Range("P32").NumberFormat = $priceFmt

$ws.Range("Q32").Value = "1:0"

# ---------------------------------------------------------------------------
# 4) Re-create the merges for the new row (A:B, C:G, H:K, L:M, N:O) to match
#    every other data row in the table.
# ---------------------------------------------------------------------------
$ws.Range("A32:B32").Merge()
$ws.Range("C32:G32").Merge()
$ws.Range("H32:K32").Merge()
$ws.Range("L32:M32").Merge()
$ws.Range("N32:O32").Merge()

# ---------------------------------------------------------------------------
# 5) Update the grand-total row (now row 36): add the new item's selling
#    price (65.00) to the previous total of 1060.37.
# ---------------------------------------------------------------------------
$ws.Range("P36").Value = 1125.37

# ---------------------------------------------------------------------------
# 6) Update the footer timestamp (now row 37): the report was regenerated a
#    few minutes later (12:26 PM -> 12:32 PM).
# ---------------------------------------------------------------------------
$ws.Range("A37").Value = "Tuesday, 29 July, 2025 12:32 PM"
